# "embed backpack system with new UI system"
#
# - Rename the UI form entry "UI_Inventory" to "UI_Backpack" (UIConfig!C9).
# - Make UIConfig the active sheet/tab (it was EntityConfig before), with its
#   selection resting on E13.
# - Leave EntityConfig's own remembered selection on C9 now that it is no
#   longer the active sheet.

$wb = $excel.ActiveWorkbook

$uiConfig = $wb.Worksheets.Item("UIConfig")
$entityConfig = $wb.Worksheets.Item("EntityConfig")

# Rename the "UI_Inventory" UI form to "UI_Backpack".
$uiConfig.Range("C9").Value = "UI_Backpack"

# Update EntityConfig's parked selection first, while it is still active.
$entityConfig.Activate()
$entityConfig.Range("C9").Select() | Out-Null

# Activate UIConfig (becomes the new active tab) and move its selection.
$uiConfig.Activate()
$uiConfig.Range("E13").Select() | Out-Null
